$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''257.49'
$ws.Range("E2").Value = '''0.95%'
$ws.Range("G2").Value = '''8'
$ws.Range("D3").Value = '''27.21'
$ws.Range("E3").Value = '''-3.22%'
$ws.Range("G3").Value = '''8'
$ws.Range("D4").Value = '''4.774'
$ws.Range("E4").Value = '''-11.45%'
$ws.Range("G4").Value = '''8'
$ws.Range("D5").Value = '''0.05962'
$ws.Range("E5").Value = '''1.94%'
$ws.Range("G5").Value = '''8'
$ws.Range("D6").Value = '''6.665'
$ws.Range("E6").Value = '''-0.64%'
$ws.Range("G6").Value = '''8'
$ws.Range("D7").Value = '''0.8696'
$ws.Range("E7").Value = '''0.35%'
$ws.Range("G7").Value = '''8'
$ws.Range("D8").Value = '''0.9607'
$ws.Range("E8").Value = '''5.63%'
$ws.Range("G8").Value = '''8'
$ws.Range("D9").Value = '''0.1407'
$ws.Range("G9").Value = '''8'
$ws.Range("D10").Value = '''0.07155'
$ws.Range("E10").Value = '''0.04%'
$ws.Range("G10").Value = '''8'
$ws.Range("D11").Value = '''0.03178'
$ws.Range("E11").Value = '''-0.09%'
$ws.Range("G11").Value = '''8'
$ws.Range("D12").Value = '''0.09253'
$ws.Range("E12").Value = '''0.29%'
$ws.Range("G12").Value = '''8'
$ws.Range("D13").Value = '''0.001538'
$ws.Range("E13").Value = '''0.07%'
$ws.Range("G13").Value = '''8'
$ws.Range("D14").Value = '''0.0006069'
$ws.Range("E14").Value = '''0.04%'
$ws.Range("G14").Value = '''8'
$ws.Range("D15").Value = '''0.006120'
$ws.Range("E15").Value = '''5.44%'
$ws.Range("G15").Value = '''8'
$ws.Range("D16").Value = '''3.480'
$ws.Range("E16").Value = '''-0.52%'
$ws.Range("G16").Value = '''8'
$ws.Range("D17").Value = '''3.186'
$ws.Range("E17").Value = '''-1.30%'
$ws.Range("G17").Value = '''8'
$ws.Range("D18").Value = '''2.219'
$ws.Range("E18").Value = '''0.77%'
$ws.Range("G18").Value = '''8'
$ws.Range("D19").Value = '''0.3131'
$ws.Range("E19").Value = '''-1.35%'
$ws.Range("G19").Value = '''8'
$ws.Range("D20").Value = '''0.03711'
$ws.Range("E20").Value = '''7.69%'
$ws.Range("G20").Value = '''8'
$ws.Range("E21").Value = '''-0.74%'
$ws.Range("G21").Value = '''8'
$ws.Range("D22").Value = '''3.811'
$ws.Range("E22").Value = '''8.26%'
$ws.Range("G22").Value = '''8'
$ws.Range("D23").Value = '''0.04225'
$ws.Range("E23").Value = '''1.80%'
$ws.Range("G23").Value = '''8'
$ws.Range("E24").Value = '''0.15%'
$ws.Range("G24").Value = '''8'
$ws.Range("D25").Value = '''0.001222'
$ws.Range("E25").Value = '''-0.05%'
$ws.Range("G25").Value = '''8'
$ws.Range("D26").Value = '''0.004502'
$ws.Range("E26").Value = '''-10.73%'
$ws.Range("G26").Value = '''8'
$ws.Range("E27").Value = '''0.05%'
$ws.Range("G27").Value = '''8'
$ws.Range("E28").Value = '''-22.93%'
$ws.Range("G28").Value = '''8'
$ws.Range("G29").Value = '''8'
$ws.Range("G30").Value = '''8'
$ws.Range("G31").Value = '''8'
$ws.Range("G32").Value = '''8'
$ws.Range("G33").Value = '''8'
$ws.Range("G34").Value = '''8'
$ws.Range("G35").Value = '''8'
$ws.Range("G36").Value = '''8'
$ws.Range("G37").Value = '''8'
$ws.Range("G38").Value = '''8'
$ws.Range("G39").Value = '''8'
$ws.Range("D40").Value = '''0.03818'
$ws.Range("E40").Value = '''-0.55%'
$ws.Range("G40").Value = '''8'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006249'
$ws.Range("E41").Value = '''9.70%'
$ws.Range("G41").Value = '''8'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1099'
$ws.Range("E42").Value = '''-0.34%'
$ws.Range("G42").Value = '''8'
$ws.Range("D43").Value = '''0.002252'
$ws.Range("E43").Value = '''-3.67%'
$ws.Range("G43").Value = '''8'
$ws.Range("D44").Value = '''0.01060'
$ws.Range("E44").Value = '''-3.22%'
$ws.Range("G44").Value = '''8'
$ws.Range("D45").Value = '''0.00005499'
$ws.Range("E45").Value = '''5.38%'
$ws.Range("G45").Value = '''8'
$ws.Range("E46").Value = '''0.05%'
$ws.Range("G46").Value = '''8'
$ws.Range("D47").Value = '''0.08850'
$ws.Range("E47").Value = '''1.15%'
$ws.Range("G47").Value = '''8'
$ws.Range("D48").Value = '''0.002363'
$ws.Range("E48").Value = '''9.63%'
$ws.Range("G48").Value = '''8'
$ws.Range("E49").Value = '''0.05%'
$ws.Range("G49").Value = '''8'
$ws.Range("E50").Value = '''0.05%'
$ws.Range("G50").Value = '''8'
$ws.Range("G51").Value = '''8'
